# Weekly update: a new market-report row for "Alcachofa" (week of 2021-11-03)
# is inserted at row 4, pushing all subsequent rows down by one (the former
# last row, 28, becomes row 29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4; Excel shifts rows 4-28
# down to 5-29 and carries their formatting with them.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44503
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112013
$ws.Range("G4").Value = "Alcachofa"
$ws.Range("H4").Value = "Madrigal"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 11500
$ws.Range("N4").Value = "$/caja 40 unidades"
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 288
$ws.Range("Q4").Value = 40
$ws.Range("R4").Value = "Hortaliza"
